$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 845, pushing existing rows 845-942 down to 846-943
$ws.Rows.Item(845).Insert()

$ws.Cells.Item(845, 1).Value = 10
$ws.Cells.Item(845, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(845, 3).Value = "La Araucanía"
$ws.Cells.Item(845, 4).Value = 45212
$ws.Cells.Item(845, 5).Value = 9
$ws.Cells.Item(845, 6).Value = 100112045
$ws.Cells.Item(845, 7).Value = "Zapallo"
$ws.Cells.Item(845, 8).Value = "Camote"
$ws.Cells.Item(845, 9).Value = "1a nueva(o)"
$ws.Cells.Item(845, 10).Value = 300
$ws.Cells.Item(845, 11).Value = 1200
$ws.Cells.Item(845, 12).Value = 1200
$ws.Cells.Item(845, 13).Value = 1200
$ws.Cells.Item(845, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(845, 15).Value = "Perú"
$ws.Cells.Item(845, 16).Value = 1200
$ws.Cells.Item(845, 17).Value = 1
$ws.Cells.Item(845, 18).Value = "Hortaliza"
